# Inserts a new weekly "Choclo" record at row 237 of Sheet1, pushing the
# existing rows 237-286 down to 238-287.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above the current row 237 (shifts 237..286 -> 238..287)
$ws.Rows.Item(237).Insert()

# Populate the newly inserted row 237 with the new record's data.
$ws.Cells.Item(237, 1).Value = 5
$ws.Cells.Item(237, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(237, 3).Value = "Maule"
$ws.Cells.Item(237, 4).Value = 44943
$ws.Cells.Item(237, 5).Value = 7
$ws.Cells.Item(237, 6).Value = 100112024
$ws.Cells.Item(237, 7).Value = "Choclo"
$ws.Cells.Item(237, 8).Value = "Choclero"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 50000
$ws.Cells.Item(237, 11).Value = 230
$ws.Cells.Item(237, 12).Value = 230
$ws.Cells.Item(237, 13).Value = 230
$ws.Cells.Item(237, 14).Value = "`$/unidad"
$ws.Cells.Item(237, 15).Value = "Región del Maule"
$ws.Cells.Item(237, 16).Value = 230
$ws.Cells.Item(237, 17).Value = 1
$ws.Cells.Item(237, 18).Value = "Hortaliza"
